# Auto-generated: apply scheduled-runner price/profit updates to the
# Faerie_Profits workbook. Each Leve table (one per crafting class sheet)
# stores scraped market prices (H:L) and computed profit (M:N) as plain
# values, so we just overwrite the affected cells per row.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8338066.5
$ws.Range("I40").Value = 50002400
$ws.Range("J40").Value = 5200
$ws.Range("K40").Value = 50002400
$ws.Range("L40").Value = 5200
$ws.Range("M40").Value = -50002225
$ws.Range("N40").Value = -5550
$ws.Range("H43").Value = 4795734
$ws.Range("I43").Value = 13165543
$ws.Range("K43").Value = 13165543
$ws.Range("M43").Value = -13165474
$ws.Range("H80").Value = 2028.8055
$ws.Range("J80").Value = 2737.1765
$ws.Range("L80").Value = 8211.529500000001
$ws.Range("N80").Value = -10207.5295
$ws.Range("H83").Value = 2028.8055
$ws.Range("J83").Value = 2737.1765
$ws.Range("L83").Value = 24634.5885
$ws.Range("N83").Value = -34618.5885
$ws.Range("H137").Value = 2388.6667
$ws.Range("I137").Value = 1212.875
$ws.Range("J137").Value = 3329.3
$ws.Range("K137").Value = 3638.625
$ws.Range("L137").Value = 9987.900000000001
$ws.Range("M137").Value = -1088.625
$ws.Range("N137").Value = -15087.9

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 8866.333000000001
$ws.Range("I5").Value = 10589.6
$ws.Range("K5").Value = 10589.6
$ws.Range("M5").Value = -10477.6
$ws.Range("H26").Value = 6751.25
$ws.Range("I26").Value = 6751.25
$ws.Range("K26").Value = 6751.25
$ws.Range("M26").Value = -6421.25
$ws.Range("H61").Value = 2758.1562
$ws.Range("I61").Value = 2282.7322
$ws.Range("K61").Value = 2282.7322
$ws.Range("M61").Value = -2070.7322
$ws.Range("H111").Value = 79999
$ws.Range("J111").Value = 79999
$ws.Range("L111").Value = 79999
$ws.Range("N111").Value = -88179
$ws.Range("H122").Value = 5332.925
$ws.Range("I122").Value = 4511.0938
$ws.Range("K122").Value = 13533.2814
$ws.Range("M122").Value = -11083.2814
$ws.Range("H125").Value = 86996.664
$ws.Range("J125").Value = 86996.664
$ws.Range("L125").Value = 86996.664
$ws.Range("N125").Value = -96836.664
$ws.Range("H136").Value = 2758.1562
$ws.Range("I136").Value = 2282.7322
$ws.Range("K136").Value = 6848.196599999999
$ws.Range("M136").Value = -4298.196599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 8866.333000000001
$ws.Range("I4").Value = 10589.6
$ws.Range("K4").Value = 10589.6
$ws.Range("M4").Value = -10474.6
$ws.Range("H22").Value = 798
$ws.Range("I22").Value = 798
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 798
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -625
$ws.Range("N22").Value = $null
$ws.Range("H94").Value = 1043.0588
$ws.Range("I94").Value = 1144.909
$ws.Range("K94").Value = 1144.909
$ws.Range("M94").Value = -693.9090000000001
$ws.Range("H107").Value = 5953.4443
$ws.Range("I107").Value = 5121.0713
$ws.Range("K107").Value = 5121.0713
$ws.Range("M107").Value = -3201.0713

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 296.65714
$ws.Range("I22").Value = 254.72728
$ws.Range("K22").Value = 254.72728
$ws.Range("M22").Value = 95.27271999999999
$ws.Range("H58").Value = 2766.9333
$ws.Range("I58").Value = 2072.6667
$ws.Range("J58").Value = 3808.3333
$ws.Range("K58").Value = 2072.6667
$ws.Range("L58").Value = 3808.3333
$ws.Range("M58").Value = -1869.6667
$ws.Range("N58").Value = -4214.3333
$ws.Range("H62").Value = 4965.8335
$ws.Range("I62").Value = 3397.7273
$ws.Range("K62").Value = 3397.7273
$ws.Range("M62").Value = -2773.7273
$ws.Range("H65").Value = 4965.8335
$ws.Range("I65").Value = 3397.7273
$ws.Range("K65").Value = 16988.6365
$ws.Range("M65").Value = -13868.6365
$ws.Range("H86").Value = 6745.5454
$ws.Range("I86").Value = 6342.6
$ws.Range("K86").Value = 6342.6
$ws.Range("M86").Value = -5219.6
$ws.Range("H89").Value = 6745.5454
$ws.Range("I89").Value = 6342.6
$ws.Range("K89").Value = 31713
$ws.Range("M89").Value = -26097
$ws.Range("H104").Value = 94264.25
$ws.Range("J104").Value = 94264.25
$ws.Range("L104").Value = 94264.25
$ws.Range("N104").Value = -99506.25
$ws.Range("H122").Value = 3346.2632
$ws.Range("J122").Value = 4753.1665
$ws.Range("L122").Value = 14259.4995
$ws.Range("N122").Value = -19159.4995
$ws.Range("H129").Value = 69888.78
$ws.Range("J129").Value = 69888.78
$ws.Range("L129").Value = 69888.78
$ws.Range("N129").Value = -79888.78
$ws.Range("H132").Value = 2399
$ws.Range("I132").Value = 2024.25
$ws.Range("K132").Value = 6072.75
$ws.Range("M132").Value = -3542.75
$ws.Range("H136").Value = 2766.9333
$ws.Range("I136").Value = 2072.6667
$ws.Range("J136").Value = 3808.3333
$ws.Range("K136").Value = 6218.000100000001
$ws.Range("L136").Value = 11424.9999
$ws.Range("M136").Value = -3668.000100000001
$ws.Range("N136").Value = -16524.9999
$ws.Range("H141").Value = 75499.57000000001
$ws.Range("J141").Value = 142833
$ws.Range("L141").Value = 142833
$ws.Range("N141").Value = -153193

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 128.0625
$ws.Range("J2").Value = 119.666664
$ws.Range("L2").Value = 717.999984
$ws.Range("N2").Value = -943.999984
$ws.Range("H38").Value = 7279
$ws.Range("J38").Value = 24993.75
$ws.Range("L38").Value = 74981.25
$ws.Range("N38").Value = -75675.25
$ws.Range("H59").Value = 2958.1667
$ws.Range("J59").Value = 7000
$ws.Range("L59").Value = 21000
$ws.Range("N59").Value = -22080

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 91720
$ws.Range("J42").Value = 91720
$ws.Range("L42").Value = 91720
$ws.Range("N42").Value = -92690
$ws.Range("H102").Value = 63707.375
$ws.Range("I102").Value = 1443.5
$ws.Range("K102").Value = 1443.5
$ws.Range("M102").Value = 178.5
$ws.Range("H113").Value = 4433.778
$ws.Range("I113").Value = 5617.6665
$ws.Range("J113").Value = 2066
$ws.Range("K113").Value = 5617.6665
$ws.Range("L113").Value = 2066
$ws.Range("M113").Value = -3447.6665
$ws.Range("N113").Value = -6406
$ws.Range("H115").Value = 91720
$ws.Range("J115").Value = 91720
$ws.Range("L115").Value = 91720
$ws.Range("N115").Value = -94070
$ws.Range("H130").Value = 81799.60000000001
$ws.Range("J130").Value = 81799.60000000001
$ws.Range("L130").Value = 81799.60000000001
$ws.Range("N130").Value = -91839.60000000001
$ws.Range("H132").Value = 5034.6953
$ws.Range("I132").Value = 4878.627
$ws.Range("K132").Value = 14635.881
$ws.Range("M132").Value = -12105.881

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 12499.297
$ws.Range("I61").Value = 9641.207
$ws.Range("K61").Value = 9641.207
$ws.Range("M61").Value = -9439.207
$ws.Range("H82").Value = 9824.8125
$ws.Range("I82").Value = 13328.25
$ws.Range("J82").Value = 6321.375
$ws.Range("K82").Value = 13328.25
$ws.Range("L82").Value = 6321.375
$ws.Range("M82").Value = -12967.25
$ws.Range("N82").Value = -7043.375
$ws.Range("H85").Value = 9824.8125
$ws.Range("I85").Value = 13328.25
$ws.Range("J85").Value = 6321.375
$ws.Range("K85").Value = 13328.25
$ws.Range("L85").Value = 6321.375
$ws.Range("M85").Value = -12080.25
$ws.Range("N85").Value = -8817.375
$ws.Range("H113").Value = 12499.297
$ws.Range("I113").Value = 9641.207
$ws.Range("K113").Value = 9641.207
$ws.Range("M113").Value = -7471.207
$ws.Range("H132").Value = 4728.5
$ws.Range("I132").Value = 3506.7273
$ws.Range("K132").Value = 10520.1819
$ws.Range("M132").Value = -7990.1819
$ws.Range("H136").Value = 4234.7383
$ws.Range("I136").Value = 4407.1714
$ws.Range("K136").Value = 13221.5142
$ws.Range("M136").Value = -10671.5142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27860.715
$ws.Range("H107").Value = 742.5185
$ws.Range("I107").Value = 558.2778
$ws.Range("K107").Value = 1674.8334
$ws.Range("M107").Value = 245.1666
$ws.Range("H125").Value = 87000
$ws.Range("J125").Value = 87000
$ws.Range("L125").Value = 87000
$ws.Range("N125").Value = -96840
$ws.Range("H126").Value = 8386.069
$ws.Range("I126").Value = 7367.84
$ws.Range("K126").Value = 22103.52
$ws.Range("M126").Value = -19633.52
$ws.Range("H132").Value = 1277.3182
$ws.Range("I132").Value = 1163.5834
$ws.Range("J132").Value = 1413.8
$ws.Range("K132").Value = 3490.7502
$ws.Range("L132").Value = 4241.4
$ws.Range("M132").Value = -960.7501999999999
$ws.Range("N132").Value = -9301.4

